$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.990.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0914"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.976"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.653.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.304.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.138.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.88%  "
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.29%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0898"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  +11.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.129"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0354"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.228"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.07%  "
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.574.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
